# Auto-generated from the commit diff: refresh the Price (D) and
# Volume(1h) (E) columns for each coin row with the new scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these columns as plain text (matching the existing sheet layout,
# where Price/Volume are stored as literal strings like "27.20" or
# "1.70%") so Excel does not silently reinterpret them as numbers or
# percentages on assignment.
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"

$ws.Range("D2").Value = "277.10"
$ws.Range("E2").Value = "0.96%"
$ws.Range("D3").Value = "27.20"
$ws.Range("E3").Value = "1.70%"
$ws.Range("D4").Value = "4.854"
$ws.Range("E4").Value = "0.06%"
$ws.Range("D5").Value = "0.06410"
$ws.Range("E5").Value = "1.36%"
$ws.Range("D6").Value = "6.932"
$ws.Range("E6").Value = "0.83%"
$ws.Range("D7").Value = "1.189"
$ws.Range("E7").Value = "-6.39%"
$ws.Range("D8").Value = "0.8758"
$ws.Range("E8").Value = "0.74%"
$ws.Range("D9").Value = "0.1540"
$ws.Range("E9").Value = "5.41%"
$ws.Range("D10").Value = "0.05129"
$ws.Range("E10").Value = "2.87%"
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").Value = "1.37%"
$ws.Range("D12").Value = "0.02953"
$ws.Range("E12").Value = "0.28%"
$ws.Range("D13").Value = "0.08976"
$ws.Range("E13").Value = "-0.64%"
$ws.Range("D14").Value = "0.001562"
$ws.Range("E14").Value = "-0.90%"
$ws.Range("D15").Value = "0.0006375"
$ws.Range("E15").Value = "0.50%"
$ws.Range("D16").Value = "0.006100"
$ws.Range("E16").Value = "1.05%"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").Value = "0.85%"
$ws.Range("D18").Value = "3.310"
$ws.Range("E18").Value = "-0.36%"
$ws.Range("E19").Value = "0.15%"
$ws.Range("E21").Value = "1.38%"
$ws.Range("D22").Value = "3.919"
$ws.Range("E22").Value = "-0.22%"
$ws.Range("D23").Value = "0.04402"
$ws.Range("E23").Value = "1.36%"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").Value = "-0.02%"
$ws.Range("D26").Value = "0.003862"
$ws.Range("E26").Value = "-9.18%"
$ws.Range("E27").Value = "8.35%"
$ws.Range("E28").Value = "15.10%"
$ws.Range("D40").Value = "0.04162"
$ws.Range("E40").Value = "2.99%"
$ws.Range("D41").Value = "0.006789"
$ws.Range("E41").Value = "1.37%"
$ws.Range("E42").Value = "0.63%"
$ws.Range("D43").Value = "0.001939"
$ws.Range("E43").Value = "-7.61%"
$ws.Range("D44").Value = "0.01150"
$ws.Range("E44").Value = "7.43%"
$ws.Range("D45").Value = "0.00005309"
$ws.Range("E45").Value = "-0.02%"
$ws.Range("E46").Value = "13.27%"
$ws.Range("D47").Value = "0.01852"
$ws.Range("E47").Value = "-7.43%"
